# Generate Report for Handoff
# Updates the localization-status report:
#  - Refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#    timestamps for the handoff batch rows (7,8,9,11,12,13).
#  - Fills in the "Priority" column ("ht") for the same handoff batch rows
#    on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)   # "Overview"
$wsZhCn     = $wb.Worksheets.Item(2)   # "zh-cn"
$wsDeDe     = $wb.Worksheets.Item(3)   # "de-de"

$rows = @(7, 8, 9, 11, 12, 13)

foreach ($r in $rows) {
    # Overview!G -> Latest HO Xliff Generate Date
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-17 00:20:50"

    # zh-cn!H -> Latest Handoff Datetime ; zh-cn!E -> Priority
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-17 00:20:45"
    $wsZhCn.Cells.Item($r, 5).Value = "ht"

    # de-de!H -> Latest Handoff Datetime ; de-de!E -> Priority
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-17 00:20:50"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}

Write-Host "Generated handoff report: timestamps and priority updated."
